# Apply the "云芽/JMG 火只木南" event insertion + refreshed "想去人数" stats
# to both the "展览" and "全部类型" worksheets (rows 15-34, 9 columns A-I).
#
# Final row data (columns B..I) for rows 1..34 after the edit.
$data = @(
        @("开始时间", "名称", "地点", "具体时间范围", "想去人数", "最低票价", "Link", "Cover"),
        @("'2024-09-22", "新余·MG动漫秋季交流会", "北湖西路158号太平洋购物东侧1-6楼 荷塘月色清华园大酒店", "2024.09.22 09:30-09.22 17:30", 19, 45, "https://show.bilibili.com/platform/detail.html?id=91781", "//i2.hdslb.com/bfs/openplatform/202409/gKdZC9dI1725323839242.jpeg"),
        @("'2024-10-01", "上饶·我和我的cos小伙伴们第二弹", "春江北大道19号时光饶城文旅小镇3号楼2层205商铺 小飞侠篮球时光校区", "2024.10.01 09:00-10.01 17:00", 89, 45, "https://show.bilibili.com/platform/detail.html?id=92049", "//i1.hdslb.com/bfs/openplatform/202409/l5VUJLIi1725545430694.jpeg"),
        @("'2024-10-01", "九江·星梦次元XACD动漫游戏博览会国庆盛典", "九瑞大道与重庆路交汇处西南角 九江国际会展中心", "2024.10.01 10:00-10.02 17:00", 266, 55, "https://show.bilibili.com/platform/detail.html?id=90732", "//i0.hdslb.com/bfs/openplatform/202408/PTmf3umB1722911829186.jpeg"),
        @("'2024-10-01", "南昌·第一届次元茂动漫展", "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆", "2024.10.01 10:00-10.01 18:00", 45, 55, "https://show.bilibili.com/platform/detail.html?id=92367", "//i0.hdslb.com/bfs/openplatform/202409/6uH7PhlF1726111899673.jpeg"),
        @("'2024-10-01", "吉安·COMIC LIFE次元假日06", "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心", "2024.10.01 09:00-10.01 18:00", 547, 55, "https://show.bilibili.com/platform/detail.html?id=90901", "//i2.hdslb.com/bfs/openplatform/202408/J892vhG21724034395965.jpeg"),
        @("'2024-10-01", "新余·逆光ZERO01动漫游戏展", "北湖中路508号 新余市会展中心", "2024.10.01 10:00-10.01 17:00", 52, 40, "https://show.bilibili.com/platform/detail.html?id=91863", "//i2.hdslb.com/bfs/openplatform/202408/AXpvlufw1724920917711.jpeg"),
        @("'2024-10-01", "景德镇·第十六届瓷都ACG动漫游戏博览会", "迎宾大道与寺山路交叉口东200米 陶博城", "2024.10.01 09:00-10.02 17:00", 2003, 55, "https://show.bilibili.com/platform/detail.html?id=91211", "//i2.hdslb.com/bfs/openplatform/202408/j9Chr8AW1723873538272.png"),
        @("'2024-10-01", "萍乡·NL15动漫游戏展", "建设西路钻石公寓西南侧60米 智博篮球馆", "2024.10.01 10:00-10.01 17:00", 66, 45, "https://show.bilibili.com/platform/detail.html?id=91666", "//i1.hdslb.com/bfs/openplatform/202409/Ltn1lbsU1725181969027.jpeg"),
        @("'2024-10-01", "赣州·十万伏特-第八届青年文化综合展览会", "东阳山路9-3号 温州街", "2024.10.01 10:00-10.02 17:00", 96, 55, "https://show.bilibili.com/platform/detail.html?id=92405", "//i1.hdslb.com/bfs/openplatform/202409/OUu2pvUR1726127605875.jpeg"),
        @("'2024-10-02", "南昌·萌卡动漫展", "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆", "2024.10.02 09:00-10.03 17:00", 4295, 65, "https://show.bilibili.com/platform/detail.html?id=89738", "//i0.hdslb.com/bfs/openplatform/202407/uqTvacSV1721621530709.jpeg"),
        @("'2024-10-02", "新余·ML动漫同好交流会10.2 ", "抱石大道561号 新余华瑞圣欧顿酒店", "2024.10.02 10:00-10.02 17:00", 33, 45, "https://show.bilibili.com/platform/detail.html?id=91340", "//i2.hdslb.com/bfs/openplatform/202408/TPp12j6F1725007454725.jpeg"),
        @("'2024-10-02", "景德镇·第十六届瓷都ACG内场—花玲&宴宁", "迎宾大道与寺山路交叉口东200米 陶博城", "2024.10.02 09:00-10.02 17:00", 280, 188, "https://show.bilibili.com/platform/detail.html?id=91080", "//i1.hdslb.com/bfs/openplatform/202408/z5XgI7ZE1724229859134.jpeg"),
        @("'2024-10-02", "鹰潭·MZD动漫游戏嘉年华", "南站路锦都金源酒店18楼 锦都金源酒店", "2024.10.02 10:00-10.02 17:00", 96, 55, "https://show.bilibili.com/platform/detail.html?id=90709", "//i0.hdslb.com/bfs/openplatform/202409/doXPtNRc1726297378846.jpeg"),
        @("'2024-10-03", " 江西·JMG（广电）第二届UP动漫游戏博览会-火只木南专场见面会", "怀玉山大道1315号 南昌绿地国际博览中心", "2024.10.03 13:00-10.03 18:00", 2, 168, "https://show.bilibili.com/platform/detail.html?id=92635", "//i2.hdslb.com/bfs/openplatform/202409/KGb5Rbt31726710245424.jpeg"),
        @("'2024-10-03", "上饶·第三届囧喵喵国风动漫展", "高铁经济试验区吴楚大道与凤凰东大道交叉口 饶派沉浸式街区(B1)", "2024.10.03 09:30-10.03 17:00", 105, 65, "https://show.bilibili.com/platform/detail.html?id=92108", "//i2.hdslb.com/bfs/openplatform/202409/hKUJo5UP1725602141242.jpeg"),
        @("'2024-10-03", "宜春·第三十七届静卿国风动漫文化汉文化展览会", "宜阳大道19号(交通银行旁) 宜春安缦文华酒店", "2024.10.03 09:00-10.03 17:00", 22, 50, "https://show.bilibili.com/platform/detail.html?id=92427", "//i2.hdslb.com/bfs/openplatform/202409/ntPXoGGS1726158288182.jpeg"),
        @("'2024-10-03", "宜春·逆光ZERO动漫游戏展", "市府北路10号  红林大酒店", "2024.10.03 10:00-10.03 17:00", 14, 40, "https://show.bilibili.com/platform/detail.html?id=91866", "//i0.hdslb.com/bfs/openplatform/202409/i1gKVM991726717574907.jpeg"),
        @("'2024-10-03", "抚州·逆光ZERO动漫游戏展", "王安石大道2466号 保利华章希尔顿逸林酒店", "2024.10.03 10:00-10.03 17:00", 63, 40, "https://show.bilibili.com/platform/detail.html?id=91865", "//i0.hdslb.com/bfs/openplatform/202409/1p0DrTb91725280390796.jpeg"),
        @("'2024-10-03", "江西·JMG（江西广电）第二届UP动漫游戏博览会", "怀玉山大道1315号 南昌绿地国际博览中心", "2024.10.03 09:00-10.05 18:00", 3061, 70, "https://show.bilibili.com/platform/detail.html?id=90599", "//i0.hdslb.com/bfs/openplatform/202408/oZpM885D1724642687206.png"),
        @("'2024-10-03", "萍乡·AU10秋至国漫展", "建设西路钻石公寓西南侧60米 智博篮球馆", "2024.10.03 10:00-10.03 17:00", 60, 45, "https://show.bilibili.com/platform/detail.html?id=92178", "//i1.hdslb.com/bfs/openplatform/202409/YWK3l0Zx1725857595232.jpeg"),
        @("'2024-10-03", "赣州·第五届半夏动漫展", "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心", "2024.10.03 10:00-10.05 17:00", 443, 49.5, "https://show.bilibili.com/platform/detail.html?id=91719", "//i2.hdslb.com/bfs/openplatform/202409/KLxI6RZQ1725270195248.jpeg"),
        @("'2024-10-03", "鹰潭·夜穹动漫游戏嘉年华", "南站路66号 回禾酒店（鹰潭火车站南站路店）", "2024.10.03 10:00-10.03 17:00", 17, 29.9, "https://show.bilibili.com/platform/detail.html?id=91960", "//i0.hdslb.com/bfs/openplatform/202409/EudXOPTz1725362358018.jpeg"),
        @("'2024-10-04", "九江·无限喵国潮动漫节", "洪垅大道 智汇欣体育中心", "2024.10.04 10:00-10.04 17:00", 15, 30, "https://show.bilibili.com/platform/detail.html?id=92370", "//i1.hdslb.com/bfs/openplatform/202409/Mebp8k9u1725460545541.jpeg"),
        @("'2024-10-04", "宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华", "宜春国际商贸城会展中心 宜春国际商贸城会展中心", "2024.10.04 10:00-10.05 17:00", 71, 55, "https://show.bilibili.com/platform/detail.html?id=91115", "//i0.hdslb.com/bfs/openplatform/202408/8TOnPvxz1723000627660.jpeg"),
        @("'2024-10-04", "江西·JMG（江西广电）第二届UP动漫游戏博览会-钱琛专场见面会", "怀玉山大道1315号 南昌绿地国际博览中心", "2024.10.04 13:00-10.04 18:00", 75, 88, "https://show.bilibili.com/platform/detail.html?id=92463", "//i1.hdslb.com/bfs/openplatform/202409/oHUl82SU1726214285169.jpeg"),
        @("'2024-10-05", "新余·星冉动漫游戏嘉年华", "长青南路 维也纳国际宴会厅（明和国际店）", "2024.10.05 10:00-10.05 17:00", 8, 29.9, "https://show.bilibili.com/platform/detail.html?id=91959", "//i0.hdslb.com/bfs/openplatform/202409/r7juTsXz1726729959726.jpeg"),
        @("'2024-10-05", "鹰潭·鹰纪元国潮动漫节", "五二五小区西侧约140米 DNA星途篮球公园", "2024.10.05 10:00-10.05 17:00", 6, 30, "https://show.bilibili.com/platform/detail.html?id=92374", "//i1.hdslb.com/bfs/openplatform/202409/0J5liYAC1725512122827.jpeg"),
        @("'2024-10-06", "南昌·星辰动漫游戏展嘉年华", "长寿路 锐成体育2＋1篮球公园", "2024.10.06 10:00-10.06 17:00", 48, 29.9, "https://show.bilibili.com/platform/detail.html?id=91961", "//i2.hdslb.com/bfs/openplatform/202409/A6PiZvto1725362505262.jpeg"),
        @("'2024-10-19", "南昌·第一届哥布林动漫游戏展——开学季&贺中秋", "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆", "2024.10.19 10:00-10.20 18:00", 194, 55, "https://show.bilibili.com/platform/detail.html?id=89240", "//i2.hdslb.com/bfs/openplatform/202409/3N3iIqVu1725270119618.jpeg"),
        @("'2024-11-16", "上饶·星河城市动漫文化节", "春江北大道时光PARK内 博悦宴会艺术中心", "2024.11.16 10:00-11.16 17:00", 8, 55, "https://show.bilibili.com/platform/detail.html?id=92572", "//i1.hdslb.com/bfs/openplatform/202409/sb00Vz8W1726637703913.jpeg"),
        @("'2024-11-16", "南昌·CM04动漫游戏博览会", "怀玉山大道1315号 南昌绿地国际博览中心", "2024.11.16 09:00-11.17 17:00", 456, 65, "https://show.bilibili.com/platform/detail.html?id=92378", "//i2.hdslb.com/bfs/openplatform/202409/N57Jfogr1725381095803.jpeg"),
        @("'2024-12-07", "南昌·云芽动漫音乐嘉年华", "怀玉山大道1315号 南昌绿地国际博览中心", "2024.12.07 09:00-12.08 18:00", 1694, 69, "https://show.bilibili.com/platform/detail.html?id=92144", "//i0.hdslb.com/bfs/openplatform/202409/2DwZA4qv1725706772865.png"),
        @("'2024-12-08", "南昌·云芽动漫音乐嘉年华·封茗囧菌内场票", "怀玉山大道1315号 南昌绿地国际博览中心", "2024.12.08 09:30-12.08 17:30", 247, 128, "https://show.bilibili.com/platform/detail.html?id=92134", "//i0.hdslb.com/bfs/openplatform/202409/eeFHJb3W1725328994111.jpeg")
    )

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new row at position 15 (shifts old rows 15-33 down to 16-34),
    # copying the formatting of row 14 onto the newly inserted row so the
    # new row's A-column keeps the bold/centered/bordered style used by
    # every other row in column A.
    $ws.Rows("15:15").Insert()
    $ws.Range("A14").Copy() | Out-Null
    $ws.Range("A15").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    # Rewrite every data row (header + 33 data rows) so values/stats are
    # fully in sync with the refreshed scrape.
    for ($i = 0; $i -lt $data.Length; $i++) {
        $row = $data[$i]
        $r = $i + 1
        $ws.Cells.Item($r, 1).Value = ($r - 1)
        for ($j = 0; $j -lt $row.Length; $j++) {
            $ws.Cells.Item($r, $j + 2).Value = $row[$j]
        }
    }
}
